$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 497.33334
$ws.Range("I18").Value = 497.33334
$ws.Range("K18").Value = 497.33334
$ws.Range("M18").Value = -213.33334
$ws.Range("H40").Value = 2133.0952
$ws.Range("I40").Value = 2163.6365
$ws.Range("J40").Value = 2099.5
$ws.Range("K40").Value = 2163.6365
$ws.Range("L40").Value = 2099.5
$ws.Range("M40").Value = -1988.6365
$ws.Range("N40").Value = -2449.5
$ws.Range("H69").Value = 6189.231
$ws.Range("I69").Value = 4666.6665
$ws.Range("K69").Value = 13999.9995
$ws.Range("M69").Value = -13125.9995
$ws.Range("H72").Value = 6189.231
$ws.Range("I72").Value = 4666.6665
$ws.Range("K72").Value = 41999.9985
$ws.Range("M72").Value = -37631.9985
$ws.Range("H76").Value = 4797.5
$ws.Range("I76").Value = 4595
$ws.Range("K76").Value = 4595
$ws.Range("M76").Value = -4280
$ws.Range("H79").Value = 4797.5
$ws.Range("I79").Value = 4595
$ws.Range("K79").Value = 4595
$ws.Range("M79").Value = -3503
$ws.Range("H138").Value = 7114.478
$ws.Range("I138").Value = 1346.9032
$ws.Range("J138").Value = 19034.133
$ws.Range("K138").Value = 4040.7096
$ws.Range("L138").Value = 57102.399
$ws.Range("M138").Value = 1099.2904
$ws.Range("N138").Value = -67382.399

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 32622.219
$ws.Range("I97").Value = 37955.555
$ws.Range("J97").Value = 3822.2
$ws.Range("K97").Value = 37955.555
$ws.Range("L97").Value = 3822.2
$ws.Range("M97").Value = -37459.555
$ws.Range("N97").Value = -4814.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 49534.19
$ws.Range("I20").Value = 54674.156
$ws.Range("J20").Value = 704.5
$ws.Range("K20").Value = 54674.156
$ws.Range("L20").Value = 704.5
$ws.Range("M20").Value = -54427.156
$ws.Range("N20").Value = -1198.5
$ws.Range("H22").Value = 426
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -277
$ws.Range("H94").Value = 648
$ws.Range("I94").Value = 526.4
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 526.4
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -75.39999999999998
$ws.Range("N94").Value = -1702
$ws.Range("H132").Value = 63338.625
$ws.Range("I132").Value = 30709
$ws.Range("J132").Value = 68000
$ws.Range("K132").Value = 30709
$ws.Range("L132").Value = 68000
$ws.Range("M132").Value = -25649
$ws.Range("N132").Value = -78120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H41").Value = 11447.143
$ws.Range("I41").Value = 4525
$ws.Range("J41").Value = 14216
$ws.Range("K41").Value = 4525
$ws.Range("L41").Value = 14216
$ws.Range("M41").Value = -4097
$ws.Range("N41").Value = -15072
$ws.Range("H50").Value = 13340
$ws.Range("J50").Value = 13340
$ws.Range("L50").Value = 13340
$ws.Range("N50").Value = -14590
$ws.Range("H51").Value = 7966.3335
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H59").Value = 21604.445
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 24180
$ws.Range("K59").Value = 1000
$ws.Range("L59").Value = 24180
$ws.Range("M59").Value = 145
$ws.Range("N59").Value = -26470
$ws.Range("H60").Value = 12154.4
$ws.Range("I60").Value = 9026
$ws.Range("J60").Value = 14240
$ws.Range("K60").Value = 9026
$ws.Range("L60").Value = 14240
$ws.Range("M60").Value = -8515
$ws.Range("N60").Value = -15262
$ws.Range("H61").Value = 7966.3335
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H68").Value = 18867.37
$ws.Range("J68").Value = 18867.37
$ws.Range("L68").Value = 18867.37
$ws.Range("N68").Value = -20365.37
$ws.Range("H71").Value = 18867.37
$ws.Range("J71").Value = 18867.37
$ws.Range("L71").Value = 56602.11
$ws.Range("N71").Value = -64090.11
$ws.Range("H74").Value = 26307.715
$ws.Range("J74").Value = 26307.715
$ws.Range("L74").Value = 26307.715
$ws.Range("N74").Value = -28055.715
$ws.Range("H77").Value = 26307.715
$ws.Range("J77").Value = 26307.715
$ws.Range("L77").Value = 78923.145
$ws.Range("N77").Value = -87659.145
$ws.Range("H107").Value = 778.26666
$ws.Range("I107").Value = 922.375
$ws.Range("J107").Value = 613.5714
$ws.Range("K107").Value = 922.375
$ws.Range("L107").Value = 613.5714
$ws.Range("M107").Value = 997.625
$ws.Range("N107").Value = -4453.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 671.64703
$ws.Range("I23").Value = 574.5
$ws.Range("J23").Value = 701.53845
$ws.Range("K23").Value = 1723.5
$ws.Range("L23").Value = 2104.61535
$ws.Range("M23").Value = -1488.5
$ws.Range("N23").Value = -2574.61535
$ws.Range("H103").Value = 2000
$ws.Range("I103").Value = 2000
$ws.Range("K103").Value = 6000
$ws.Range("M103").Value = -5121
$ws.Range("H113").Value = 639.80554
$ws.Range("I113").Value = 580.8333
$ws.Range("J113").Value = 669.2917
$ws.Range("K113").Value = 1742.4999
$ws.Range("L113").Value = 2007.8751
$ws.Range("M113").Value = 427.5001
$ws.Range("N113").Value = -6347.8751
$ws.Range("H123").Value = 4340.125
$ws.Range("J123").Value = 4948.5
$ws.Range("L123").Value = 14845.5
$ws.Range("N123").Value = -19745.5
$ws.Range("H129").Value = 15941384
$ws.Range("J129").Value = 384989.53
$ws.Range("L129").Value = 1154968.59
$ws.Range("N129").Value = -1164968.59
$ws.Range("H131").Value = 863.95
$ws.Range("I131").Value = 443.33334
$ws.Range("J131").Value = 876.95874
$ws.Range("K131").Value = 1330.00002
$ws.Range("L131").Value = 2630.87622
$ws.Range("M131").Value = 3709.99998
$ws.Range("N131").Value = -12710.87622
$ws.Range("H137").Value = 17340868
$ws.Range("I137").Value = 4524.0835
$ws.Range("J137").Value = 33343646
$ws.Range("K137").Value = 13572.2505
$ws.Range("L137").Value = 100030938
$ws.Range("M137").Value = -8472.250499999998
$ws.Range("N137").Value = -100041138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 14056.75
$ws.Range("J136").Value = 14056.75
$ws.Range("L136").Value = 42170.25
$ws.Range("N136").Value = -47270.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4899.143
$ws.Range("I7").Value = 2248.5
$ws.Range("J7").Value = 8433.333000000001
$ws.Range("K7").Value = 2248.5
$ws.Range("L7").Value = 8433.333000000001
$ws.Range("M7").Value = -2136.5
$ws.Range("N7").Value = -8657.333000000001
$ws.Range("H68").Value = 2837.9412
$ws.Range("J68").Value = 3558.6365
$ws.Range("L68").Value = 3558.6365
$ws.Range("N68").Value = -5056.636500000001
$ws.Range("H71").Value = 2837.9412
$ws.Range("J71").Value = 3558.6365
$ws.Range("L71").Value = 17793.1825
$ws.Range("N71").Value = -25281.1825
$ws.Range("H126").Value = 4899.143
$ws.Range("I126").Value = 2248.5
$ws.Range("J126").Value = 8433.333000000001
$ws.Range("K126").Value = 6745.5
$ws.Range("L126").Value = 25299.999
$ws.Range("M126").Value = -4275.5
$ws.Range("N126").Value = -30239.999
